$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 18, shifting the existing rows 18-20 down to 19-21.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new record's data.
$ws.Range("A18").Value = 12
$ws.Range("B18").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C18").Value = "Metropolitana"
$ws.Range("D18").Value2 = 44504
$ws.Range("E18").Value = 13
$ws.Range("F18").Value = 100112028
$ws.Range("G18").Value = "Sandia"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 200
$ws.Range("K18").Value = 800
$ws.Range("L18").Value = 800
$ws.Range("M18").Value = 800
$ws.Range("N18").Value = "$/kilo (volumen en unidades)"
$ws.Range("O18").Value = "Perú"
$ws.Range("P18").Value = 800
$ws.Range("Q18").Value = 1
$ws.Range("R18").Value = "Hortaliza"
